$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-23 (columns A=count, B=bot, C=edit_ts)
$data = @(
    @("2",   $false, "15:10"),
    @("465", $false, "15:11"),
    @("434", $false, "15:12"),
    @("416", $false, "15:13"),
    @("422", $false, "15:14"),
    @("374", $false, "15:15"),
    @("867", $false, "15:16"),
    @("473", $false, "15:17"),
    @("726", $false, "15:18"),
    @("945", $false, "15:19"),
    @("376", $false, "15:20"),
    @("14",  $true,  "15:10"),
    @("486", $true,  "15:11"),
    @("502", $true,  "15:12"),
    @("485", $true,  "15:13"),
    @("470", $true,  "15:14"),
    @("501", $true,  "15:15"),
    @("473", $true,  "15:16"),
    @("494", $true,  "15:17"),
    @("450", $true,  "15:18"),
    @("441", $true,  "15:19"),
    @("184", $true,  "15:20")
)

$lastRow = 1 + $data.Count

# Make sure column A (count) is stored as text so numeric-looking values
# like "2" or "465" are not reinterpreted as numbers.
$ws.Range("A2:A$lastRow").NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("A1:C$lastRow").Select()
